# Auto-generated edit script: updates currentAveragePrice / Leve price / profit columns
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H9").Value = 788.93335
$ws.Range("I9").Value = 859.9
$ws.Range("J9").Value = 647
$ws.Range("K9").Value = 859.9
$ws.Range("L9").Value = 647
$ws.Range("M9").Value = -690.9
$ws.Range("N9").Value = -985
$ws.Range("H16").Value = 9
$ws.Range("I16").Value = 9
$ws.Range("K16").Value = 9
$ws.Range("M16").Value = 221
$ws.Range("H18").Value = 7500
$ws.Range("I18").Value = 7500
$ws.Range("K18").Value = 7500
$ws.Range("M18").Value = -7216
$ws.Range("H51").Value = 3583
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 2599.6
$ws.Range("I62").Value = 2233
$ws.Range("J62").Value = 3149.5
$ws.Range("K62").Value = 2233
$ws.Range("L62").Value = 3149.5
$ws.Range("M62").Value = -1609
$ws.Range("N62").Value = -4397.5
$ws.Range("H65").Value = 2599.6
$ws.Range("I65").Value = 2233
$ws.Range("J65").Value = 3149.5
$ws.Range("K65").Value = 11165
$ws.Range("L65").Value = 15747.5
$ws.Range("M65").Value = -8045
$ws.Range("N65").Value = -21987.5
$ws.Range("H74").Value = 4164.846
$ws.Range("I74").Value = 3678.5833
$ws.Range("K74").Value = 3678.5833
$ws.Range("M74").Value = -2742.5833
$ws.Range("H77").Value = 4164.846
$ws.Range("I77").Value = 3678.5833
$ws.Range("K77").Value = 18392.9165
$ws.Range("M77").Value = -13712.9165
$ws.Range("H135").Value = 475.85715
$ws.Range("I135").Value = 475.85715
$ws.Range("K135").Value = 4282.71435
$ws.Range("M135").Value = -1747.71435
$ws.Range("H137").Value = 1382.9231
$ws.Range("I137").Value = 1535.9
$ws.Range("J137").Value = 873
$ws.Range("K137").Value = 4607.700000000001
$ws.Range("L137").Value = 2619
$ws.Range("M137").Value = -2057.700000000001
$ws.Range("N137").Value = -7719

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1365.7059
$ws.Range("I2").Value = 1365.7059
$ws.Range("K2").Value = 1365.7059
$ws.Range("M2").Value = -1252.7059
$ws.Range("H5").Value = 180.16667
$ws.Range("J5").Value = 240
$ws.Range("L5").Value = 240
$ws.Range("N5").Value = -464
$ws.Range("H61").Value = 2077.9
$ws.Range("I61").Value = 2097.375
$ws.Range("K61").Value = 2097.375
$ws.Range("M61").Value = -1885.375
$ws.Range("H116").Value = 1365.7059
$ws.Range("I116").Value = 1365.7059
$ws.Range("K116").Value = 1365.7059
$ws.Range("M116").Value = 928.2941000000001
$ws.Range("H136").Value = 2077.9
$ws.Range("I136").Value = 2097.375
$ws.Range("K136").Value = 6292.125
$ws.Range("M136").Value = -3742.125

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1365.7059
$ws.Range("I3").Value = 1365.7059
$ws.Range("K3").Value = 1365.7059
$ws.Range("M3").Value = -1251.7059
$ws.Range("H4").Value = 180.16667
$ws.Range("J4").Value = 240
$ws.Range("L4").Value = 240
$ws.Range("N4").Value = -470
$ws.Range("H5").Value = 453
$ws.Range("I5").Value = 901
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 901
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = -788
$ws.Range("N5").Value = -231
$ws.Range("H7").Value = 50
$ws.Range("J7").Value = 50
$ws.Range("L7").Value = 50
$ws.Range("N7").Value = -276
$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 1500
$ws.Range("K21").Value = 1500
$ws.Range("M21").Value = -1265
$ws.Range("H107").Value = 936.6667
$ws.Range("I107").Value = 936.6667
$ws.Range("K107").Value = 936.6667
$ws.Range("M107").Value = 983.3333
$ws.Range("H134").Value = 2811.125
$ws.Range("I134").Value = 3097.8
$ws.Range("J134").Value = 2333.3333
$ws.Range("K134").Value = 9293.400000000001
$ws.Range("L134").Value = 6999.999899999999
$ws.Range("M134").Value = -6758.400000000001
$ws.Range("N134").Value = -12069.9999

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H58").Value = 5332.8335
$ws.Range("I58").Value = 5499.5
$ws.Range("K58").Value = 5499.5
$ws.Range("M58").Value = -5296.5
$ws.Range("H99").Value = 2999.8
$ws.Range("I99").Value = 2999.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2999.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1501.8
$ws.Range("N99").ClearContents()
$ws.Range("H108").Value = 18000
$ws.Range("I108").Value = 5000
$ws.Range("K108").Value = 5000
$ws.Range("M108").Value = -1160
$ws.Range("H126").Value = 2999.8
$ws.Range("I126").Value = 2999.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8999.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6529.400000000001
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 5332.8335
$ws.Range("I136").Value = 5499.5
$ws.Range("K136").Value = 16498.5
$ws.Range("M136").Value = -13948.5

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H7").Value = 13040.4375
$ws.Range("I7").Value = 15876.692
$ws.Range("J7").Value = 750
$ws.Range("K7").Value = 47630.076
$ws.Range("L7").Value = 2250
$ws.Range("M7").Value = -47518.076
$ws.Range("N7").Value = -2474
$ws.Range("H9").Value = 1115.6
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 1240.4546
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 3721.3638
$ws.Range("M9").Value = -376
$ws.Range("N9").Value = -4169.3638
$ws.Range("H12").Value = 969.6842
$ws.Range("I12").Value = 139.22223
$ws.Range("J12").Value = 1717.1
$ws.Range("K12").Value = 417.66669
$ws.Range("L12").Value = 5151.299999999999
$ws.Range("M12").Value = -244.66669
$ws.Range("N12").Value = -5497.299999999999
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 15000
$ws.Range("M31").Value = -14712
$ws.Range("H34").Value = 3062.077
$ws.Range("I34").Value = 218.16667
$ws.Range("K34").Value = 654.50001
$ws.Range("M34").Value = -570.50001
$ws.Range("H81").Value = 11239
$ws.Range("I81").Value = 1695
$ws.Range("K81").Value = 5085
$ws.Range("M81").Value = -3962
$ws.Range("H84").Value = 11239
$ws.Range("I84").Value = 1695
$ws.Range("K84").Value = 15255
$ws.Range("M84").Value = -9639

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H3").Value = 17851672
$ws.Range("I3").Value = 35627750
$ws.Range("K3").Value = 35627750
$ws.Range("M3").Value = -35627634
$ws.Range("H4").Value = 1003
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 4800
$ws.Range("I5").Value = 4800
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4800
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4688
$ws.Range("N5").ClearContents()
$ws.Range("H113").Value = 2994.625
$ws.Range("I113").Value = 2995.4
$ws.Range("J113").Value = 2993.3333
$ws.Range("K113").Value = 2995.4
$ws.Range("L113").Value = 2993.3333
$ws.Range("M113").Value = -825.4000000000001
$ws.Range("N113").Value = -7333.3333
$ws.Range("H126").Value = 1978
$ws.Range("I126").Value = 1978
$ws.Range("K126").Value = 5934
$ws.Range("M126").Value = -3464

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 4710.357
$ws.Range("I7").Value = 4803.5386
$ws.Range("K7").Value = 4803.5386
$ws.Range("M7").Value = -4691.5386
$ws.Range("H22").Value = 1248.5
$ws.Range("I22").Value = 1228.2
$ws.Range("K22").Value = 1228.2
$ws.Range("M22").Value = -933.2
$ws.Range("H27").Value = 1248.5
$ws.Range("I27").Value = 1228.2
$ws.Range("K27").Value = 1228.2
$ws.Range("M27").Value = -1121.2
$ws.Range("H30").Value = 826.5
$ws.Range("I30").Value = 894.5714
$ws.Range("K30").Value = 894.5714
$ws.Range("M30").Value = -786.5714
$ws.Range("H40").Value = 1369.4166
$ws.Range("I40").Value = 1369.4166
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1369.4166
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1233.4166
$ws.Range("N40").ClearContents()
$ws.Range("H126").Value = 4710.357
$ws.Range("I126").Value = 4803.5386
$ws.Range("K126").Value = 14410.6158
$ws.Range("M126").Value = -11940.6158
$ws.Range("H132").Value = 4024.6667
$ws.Range("I132").Value = 4049.8
$ws.Range("J132").Value = 3899
$ws.Range("K132").Value = 12149.4
$ws.Range("L132").Value = 11697
$ws.Range("M132").Value = -9619.400000000001
$ws.Range("N132").Value = -16757

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H34").Value = 10012.5
$ws.Range("I34").Value = 10012.5
$ws.Range("K34").Value = 10012.5
$ws.Range("M34").Value = -9809.5
$ws.Range("H132").Value = 2452.4736
$ws.Range("I132").Value = 2705.1765
$ws.Range("K132").Value = 8115.529500000001
$ws.Range("M132").Value = -5585.529500000001
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H136").Value = 3984.1162
$ws.Range("I136").Value = 3778.2856
$ws.Range("J136").Value = 4180.591
$ws.Range("K136").Value = 11334.8568
$ws.Range("L136").Value = 12541.773
$ws.Range("M136").Value = -8784.856800000001
$ws.Range("N136").Value = -17641.773
